# Hoàn thiện Ngoại Trú
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsCheck = $wb.Worksheets.Item("Check")

# Update Data sheet row 2
$wsData.Range("A2").Value = 3017
$wsData.Range("E2").Value = 46200608017
$wsData.Range("X2").Value = "DN4127460130017"

# Update Check sheet row 2
$wsCheck.Range("A2").Value = 3017
$wsCheck.Range("C2").Value = "DN4127460130017"
